$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3870
$ws.Range("C2").Value = 13
$ws.Range("D2").Value = 144
$ws.Range("E2").Value = 354
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 842
$ws.Range("I2").Value = 21222
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 181
$ws.Range("L2").Value = 19073

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 179
$ws.Range("E3").Value = 182
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 55
$ws.Range("I3").Value = 8315
$ws.Range("J3").Value = 3

$ws.Range("B4").Value = 267

$ws.Range("B5").Value = 24
$ws.Range("D5").Value = 108
$ws.Range("E5").Value = 108
$ws.Range("I5").Value = 9167
$ws.Range("J5").Value = 2

$ws.Range("B6").Value = 268
$ws.Range("D6").Value = 31
$ws.Range("E6").Value = 32
$ws.Range("F6").Value = 1
$ws.Range("H6").Value = 313
$ws.Range("I6").Value = 3871
$ws.Range("J6").Value = 0

$ws.Range("B7").Value = 264
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("I7").Value = 0

$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 83
$ws.Range("E8").Value = 95
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 638
$ws.Range("I8").Value = 3735

$ws.Range("B9").Value = 776
$ws.Range("D9").Value = 13
$ws.Range("E9").Value = 13
$ws.Range("I9").Value = 10000

$ws.Range("B10").Value = 773
$ws.Range("D10").Value = 33
$ws.Range("E10").Value = 36
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1458
$ws.Range("I10").Value = 19448
$ws.Range("J10").Value = 3

$ws.Range("B11").Value = 962
$ws.Range("D11").Value = 111
$ws.Range("E11").Value = 116
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1423
$ws.Range("I11").Value = 33884
$ws.Range("J11").Value = 1

$ws.Range("B12").Value = 840
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 14
$ws.Range("E12").Value = 14
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 8000
$ws.Range("J12").Value = 0

$ws.Range("B13").Value = 251
$ws.Range("C13").Value = 25
$ws.Range("D13").Value = 129
$ws.Range("E13").Value = 195
$ws.Range("F13").Value = 39
$ws.Range("H13").Value = 3952
$ws.Range("I13").Value = 3839
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0

$ws.Range("D14").Value = 68
$ws.Range("E14").Value = 68
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 9167
$ws.Range("J14").Value = 4

$ws.Range("B15").Value = 266
$ws.Range("D15").Value = 13
$ws.Range("E15").Value = 13
$ws.Range("I15").Value = 3077

$ws.Range("B16").Value = 22
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 134
$ws.Range("E16").Value = 190
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 4035
$ws.Range("I16").Value = 14604
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 32
$ws.Range("L16").Value = 2370

$ws.Range("B17").Value = 65
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = 23
$ws.Range("I17").Value = 5217

$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 213
$ws.Range("E18").Value = 243
$ws.Range("F18").Value = 5
$ws.Range("H18").Value = 207
$ws.Range("I18").Value = 5822
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 18
$ws.Range("L18").Value = 741

$ws.Range("B19").Value = 782
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 76
$ws.Range("E19").Value = 80
$ws.Range("F19").Value = 4
$ws.Range("H19").Value = 563
$ws.Range("I19").Value = 6915

$ws.Range("B20").Value = 1401
$ws.Range("D20").Value = 97
$ws.Range("E20").Value = 124
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 2137
$ws.Range("I20").Value = 13261
$ws.Range("J20").Value = 2

$ws.Range("B21").Value = 483
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 18
$ws.Range("I21").Value = 28667
$ws.Range("J21").Value = 1

$ws.Range("B22").Value = 190
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 0

$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 3
$ws.Range("I23").Value = 10000
$ws.Range("J23").Value = 0

$ws.Range("B24").Value = 995
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 10
$ws.Range("I24").Value = 18750

$ws.Range("B25").Value = 76
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = 29
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 5385
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1034

$ws.Range("B26").Value = 135
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 30
$ws.Range("I26").Value = 4286
$ws.Range("J26").Value = 0

$ws.Range("B27").Value = 13
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = 14
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 1429
$ws.Range("I27").Value = 3000
